# Applies the cell-value updates described in the commit diff for
# Jogos_do_Dia_Betfair_Back_Lay_2025-11-19.xlsx (odds refresh on Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2" = 3.35
    "G2" = 3.8
    "H2" = 2.36
    "I2" = 2.64
    "J2" = 2.94
    "N2" = 2.78
    "P2" = 1.55
    "Q2" = 2.56
    "R2" = 1.2
    "V2" = 1.61
    "W2" = 1.37
    "X2" = 9.6
    "Y2" = 8.199999999999999
    "Z2" = 15.5
    "AB2" = 11
    "AD2" = 12.5
    "AE2" = 130
    "AF2" = 1000
    "AG2" = 16.5
    "AH2" = 40
    "AI2" = 190
    "AJ2" = 900
    "AK2" = 220
    "AN2" = 1000
    "F3" = 2.64
    "G3" = 2.92
    "H3" = 2.98
    "I3" = 3.1
    "J3" = 3.05
    "L3" = 1.49
    "M3" = 1.09
    "N3" = 3.1
    "O3" = 1.41
    "P3" = 1.69
    "Q3" = 2.26
    "R3" = 1.25
    "S3" = 4.3
    "T3" = 1.87
    "U3" = 1.95
    "V3" = 1.47
    "X3" = 11.5
    "Z3" = 19.5
    "AA3" = 400
    "AB3" = 9.800000000000001
    "AD3" = 13.5
    "AE3" = 280
    "AF3" = 18.5
    "AH3" = 19.5
    "AK3" = 140
    "AL3" = 260
    "AN3" = 600
    "AO3" = 600
    "J4" = 5.9
    "K4" = 6
    "L4" = 1.34
    "N4" = 4.7
    "P4" = 2.26
    "Q4" = 1.76
    "R4" = 1.49
    "S4" = 2.96
    "T4" = 2.28
    "U4" = 1.73
    "W4" = 4
    "X4" = 20
    "Y4" = 40
    "Z4" = 130
    "AA4" = 580
    "AB4" = 8
    "AC4" = 13.5
    "AD4" = 46
    "AE4" = 200
    "AF4" = 7.4
    "AH4" = 34
    "AI4" = 200
    "AJ4" = 9.800000000000001
    "AM4" = 250
    "AN4" = 5.2
    "AO4" = 370
    "F5" = 1.87
    "G5" = 1.93
    "H5" = 5.2
    "I5" = 5.5
    "J5" = 3.4
    "K5" = 3.6
    "P5" = 1.59
    "Q5" = 2.58
    "G6" = 2.94
    "I6" = 3.15
    "J6" = 3.1
    "K6" = 3.55
    "N6" = 3.45
    "O6" = 1.32
    "R6" = 1.34
    "S6" = 3.35
    "T6" = 1.71
    "V6" = 1.47
    "W6" = 1.52
    "F7" = 1.96
    "G7" = 2.04
    "I7" = 5.4
    "J7" = 3.25
    "Q7" = 2.58
    "R7" = 1.19
    "U7" = 1.71
    "V7" = 1.23
    "W7" = 1.96
    "AC7" = 14
    "AD7" = 24
    "AG7" = 30
    "AK7" = 75
    "F8" = 4.3
    "G8" = 4.5
    "H8" = 2.12
    "I8" = 2.14
    "J8" = 3.3
    "K8" = 3.35
    "L8" = 1.58
    "N8" = 2.78
    "P8" = 1.58
    "R8" = 1.2
    "S8" = 5.5
    "U8" = 1.76
    "V8" = 1.87
    "W8" = 1.29
    "Z8" = 11
    "AA8" = 26
    "AE8" = 29
    "AF8" = 29
    "AG8" = 18.5
    "AJ8" = 110
    "AK8" = 75
    "AL8" = 100
    "AM8" = 200
    "AO8" = 26
    "F9" = 2.12
    "G9" = 2.14
    "H9" = 3.9
    "I9" = 3.95
    "J9" = 3.6
    "K9" = 3.65
    "N9" = 3.75
    "O9" = 1.34
    "P9" = 1.91
    "Q9" = 2.04
    "S9" = 3.7
    "T9" = 1.85
    "V9" = 1.33
    "W9" = 1.87
    "X9" = 14
    "Y9" = 14
    "Z9" = 27
    "AA9" = 80
    "AB9" = 9.199999999999999
    "AC9" = 7.8
    "AD9" = 15.5
    "AE9" = 50
    "AF9" = 12.5
    "AG9" = 10.5
    "AH9" = 18.5
    "AI9" = 60
    "AJ9" = 25
    "AK9" = 22
    "AL9" = 38
    "AM9" = 110
    "AN9" = 17
    "AO9" = 48
    "F10" = 2.28
    "G10" = 2.3
    "H10" = 3.55
    "I10" = 3.65
    "K10" = 3.55
    "O10" = 1.35
    "P10" = 1.91
    "R10" = 1.34
    "T10" = 1.85
    "V10" = 1.38
    "W10" = 1.77
    "Y10" = 13.5
    "AB10" = 9.6
    "AD10" = 15.5
    "AE10" = 44
    "AF10" = 13.5
    "AJ10" = 29
    "AK10" = 24
    "AN10" = 19
    "AO10" = 44
    "F11" = 2.58
    "G11" = 2.76
    "H11" = 2.88
    "I11" = 3.15
    "J11" = 3.3
    "K11" = 3.5
    "L11" = 1.4
    "M11" = 1.07
    "P11" = 1.98
    "Q11" = 1.93
    "S11" = 3.3
    "T11" = 1.7
    "V11" = 1.47
    "W11" = 1.57
    "AC11" = 1000
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
